$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# Update the appearance of the vj3vs28 question (row 23) from "minimal" to "likert"
$survey.Range("G23").Value = "likert"

# Update the vj3vs28 choice list labels to alternate/compact Likert formats
$choices.Range("C45").Value = "-"
$choices.Range("C46").Value = "Neutral"
$choices.Range("C47").Value = "-"
